# Updated cryptos list on Wed Sep 25 21:36:52 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a price cell (column D) as exact text, avoiding Excel's
# automatic conversion of numeric-looking strings into floating point
# numbers (which would corrupt the exact decimal representation).
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

# Row 2 - Bitcoin
Set-TextValue "D2" "63.492.98"
$ws.Range("E2").Value = "  -0.96%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.583.52"
$ws.Range("E3").Value = "  -2.46%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.04%  "

# Row 5 - BNB
Set-TextValue "D5" "589.76"
$ws.Range("E5").Value = "  -2.91%  "

# Row 6 - Solana
Set-TextValue "D6" "150.24"
$ws.Range("E6").Value = "  +1.24%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.01%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -0.50%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +1.42%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  +2.07%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  -1.12%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  -0.43%  "

# Row 13 - Avalanche
Set-TextValue "D13" "27.60"
$ws.Range("E13").Value = "  +0.16%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue "D14" "3.049.42"
$ws.Range("E14").Value = "  -2.43%  "

# Row 15 - WrappedBTC
Set-TextValue "D15" "63.344.57"
$ws.Range("E15").Value = "  -0.98%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  +5.77%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "2.597.85"
$ws.Range("E17").Value = "  -1.69%  "

# Row 18 - Chainlink
Set-TextValue "D18" "12.23"
$ws.Range("E18").Value = "  +2.29%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  +3.02%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "346.86"
$ws.Range("E20").Value = "  +0.07%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  -0.68%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  -0.05%  "

# Row 23 - Litecoin
Set-TextValue "D23" "67.30"
$ws.Range("E23").Value = "  +1.45%  "

# Row 24 - SuiNetwork
Set-TextValue "D24" "1.69"
$ws.Range("E24").Value = "  +0.88%  "

# Row 25 - InternetComputer(DFINITY)
$ws.Range("E25").Value = "  -1.98%  "

# Row 26 - Fetch.AI
$ws.Range("E26").Value = "  -2.99%  "

# Row 27 - Bittensor
Set-TextValue "D27" "553.22"
$ws.Range("E27").Value = "  -0.58%  "

# Row 28 - now Kaspa (was Aptos)
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D28" "0.163"
$ws.Range("E28").Value = "  +0.65%  "

# Row 29 - now Aptos (was Kaspa)
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D29" "8.06"
$ws.Range("E29").Value = "  -1.05%  "

# Row 30 - Binance-PegBSC-USD
$ws.Range("E30").Value = "  -0.09%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -1.36%  "

# Row 32 - PEPE
$ws.Range("E32").Value = "  +1.40%  "

# Row 33 - ImmutableX
$ws.Range("E33").Value = "  -0.66%  "

# Row 34 - NEARProtocol
Set-TextValue "D34" "5.24"
$ws.Range("E34").Value = "  -1.77%  "

# Row 35 - Monero
Set-TextValue "D35" "166.51"
$ws.Range("E35").Value = "  -1.09%  "

# Row 36 - PolygonEcosystemToken
$ws.Range("E36").Value = "  +1.70%  "

# Row 37 - FirstDigitalUSD
$ws.Range("E37").Value = "  -0.18%  "

# Row 38 - EthereumClassic
Set-TextValue "D38" "19.54"
$ws.Range("E38").Value = "  +1.11%  "

# Row 39 - Stacks
$ws.Range("E39").Value = "  -1.81%  "

# Row 40 - USDe
$ws.Range("E40").Value = "  -0.02%  "

# Row 41 - Aave
Set-TextValue "D41" "166.25"
$ws.Range("E41").Value = "  -0.44%  "

# Row 42 - OKB
Set-TextValue "D42" "39.84"
$ws.Range("E42").Value = "  -1.32%  "

# Row 43 - Filecoin
Set-TextValue "D43" "4.03"
$ws.Range("E43").Value = "  +5.01%  "

# Row 44 - InjectiveProtocol
Set-TextValue "D44" "23.08"
$ws.Range("E44").Value = "  +4.58%  "

# Row 45 - Hedera
Set-TextValue "D45" "0.0593"
$ws.Range("E45").Value = "  +3.84%  "

# Row 46 - dogwifhat
Set-TextValue "D46" "2.11"
$ws.Range("E46").Value = "  +6.09%  "

# Row 47 - Mantle
$ws.Range("E47").Value = "  +0.09%  "

# Row 48 - VeChain
Set-TextValue "D48" "0.0252"
$ws.Range("E48").Value = "  +2.71%  "

# Row 49 - Stellar
Set-TextValue "D49" "0.0964"
$ws.Range("E49").Value = "  +0.27%  "

# Row 50 - EnergySwap
Set-TextValue "D50" "19.17"
$ws.Range("E50").Value = "  +0.65%  "

# Row 51 - BabyDogeCoin
$ws.Range("D51").Value = "0.0₆0232"
$ws.Range("E51").Value = "  +15.45%  "
